# Auto-generated Excel COM-interop script
# Updates market-price snapshot values (columns H-N) on several leve-profit
# sheets, matching the scheduled-runner commit described in the diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 14468.947
$ws.Range("I40").Value = 4558.5713
$ws.Range("J40").Value = 20250
$ws.Range("K40").Value = 4558.5713
$ws.Range("L40").Value = 20250
$ws.Range("M40").Value = -4383.5713
$ws.Range("N40").Value = -20600
# Row 41
$ws.Range("H41").Value = 498
$ws.Range("I41").Value = 262.3
$ws.Range("J41").Value = 1087.25
$ws.Range("K41").Value = 262.3
$ws.Range("L41").Value = 1087.25
$ws.Range("M41").Value = 177.7
$ws.Range("N41").Value = -1967.25
# Row 80
$ws.Range("H80").Value = 488.25
$ws.Range("I80").Value = 597.5
$ws.Range("J80").Value = 379
$ws.Range("K80").Value = 1792.5
$ws.Range("L80").Value = 1137
$ws.Range("M80").Value = -794.5
$ws.Range("N80").Value = -3133
# Row 83
$ws.Range("H83").Value = 488.25
$ws.Range("I83").Value = 597.5
$ws.Range("J83").Value = 379
$ws.Range("K83").Value = 5377.5
$ws.Range("L83").Value = 3411
$ws.Range("M83").Value = -385.5
$ws.Range("N83").Value = -13395
# Row 94
$ws.Range("H94").Value = 566.5
$ws.Range("I94").Value = 589.8
$ws.Range("J94").Value = 450
$ws.Range("K94").Value = 589.8
$ws.Range("L94").Value = 450
$ws.Range("M94").Value = -138.8
$ws.Range("N94").Value = -1352
# Row 96
$ws.Range("H96").Value = 1301.7
$ws.Range("I96").Value = 897.6
$ws.Range("K96").Value = 2692.8
$ws.Range("M96").Value = -1319.8
# Row 97
$ws.Range("H97").Value = 999.5
$ws.Range("J97").Value = 999.5
$ws.Range("L97").Value = 2998.5
$ws.Range("N97").Value = -3990.5
# Row 113
$ws.Range("H113").Value = 5004
$ws.Range("J113").Value = 5177.4
$ws.Range("L113").Value = 5177.4
$ws.Range("N113").Value = -11685.4
# Row 116
$ws.Range("H116").Value = 4871.4287
$ws.Range("I116").Value = 4898.3335
$ws.Range("J116").Value = 4851.25
$ws.Range("K116").Value = 4898.3335
$ws.Range("L116").Value = 4851.25
$ws.Range("M116").Value = -1456.3335
$ws.Range("N116").Value = -11735.25
# Row 121
$ws.Range("H121").Value = 3239.818
$ws.Range("J121").Value = 3239.818
$ws.Range("L121").Value = 9719.454000000002
$ws.Range("N121").Value = -13213.454
# Row 132
$ws.Range("H132").Value = 3706.8518
$ws.Range("I132").Value = 3703.6924
$ws.Range("K132").Value = 11111.0772
$ws.Range("M132").Value = -8581.0772
# Row 137
$ws.Range("H137").Value = 1763.9231
$ws.Range("I137").Value = 1325.8334
$ws.Range("K137").Value = 3977.5002
$ws.Range("M137").Value = -1427.5002
# Row 138
$ws.Range("H138").Value = 3495.415
$ws.Range("J138").Value = 4263.4
$ws.Range("L138").Value = 12790.2
$ws.Range("N138").Value = -23070.2

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6276.8887
$ws.Range("I32").Value = 633.75
$ws.Range("K32").Value = 633.75
$ws.Range("M32").Value = -346.75
# Row 45
$ws.Range("H45").Value = 14999.333
$ws.Range("I45").Value = 22448.8
$ws.Range("K45").Value = 22448.8
$ws.Range("M45").Value = -22071.8
# Row 61
$ws.Range("H61").Value = 6276.607
$ws.Range("J61").Value = 7498.6
$ws.Range("L61").Value = 7498.6
$ws.Range("N61").Value = -7922.6
# Row 97
$ws.Range("H97").Value = 1727.5
$ws.Range("I97").Value = 638.1667
$ws.Range("J97").Value = 3688.3
$ws.Range("K97").Value = 638.1667
$ws.Range("L97").Value = 3688.3
$ws.Range("M97").Value = -142.1667
$ws.Range("N97").Value = -4680.3
# Row 132
$ws.Range("H132").Value = 5767.4116
$ws.Range("I132").Value = 5673.375
$ws.Range("K132").Value = 17020.125
$ws.Range("M132").Value = -14490.125
# Row 136
$ws.Range("H136").Value = 6276.607
$ws.Range("J136").Value = 7498.6
$ws.Range("L136").Value = 22495.8
$ws.Range("N136").Value = -27595.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 28775.6
$ws.Range("I16").Value = 969.75
$ws.Range("J16").Value = 139999
$ws.Range("K16").Value = 969.75
$ws.Range("L16").Value = 139999
$ws.Range("M16").Value = -682.75
$ws.Range("N16").Value = -140573
# Row 94
$ws.Range("H94").Value = 1067.3334
$ws.Range("J94").Value = 1368.6666
$ws.Range("L94").Value = 1368.6666
$ws.Range("N94").Value = -2270.6666
# Row 113
$ws.Range("H113").Value = 28775.6
$ws.Range("I113").Value = 969.75
$ws.Range("J113").Value = 139999
$ws.Range("K113").Value = 969.75
$ws.Range("L113").Value = 139999
$ws.Range("M113").Value = 1200.25
$ws.Range("N113").Value = -144339
# Row 123
$ws.Range("H123").Value = 89999
$ws.Range("J123").Value = 89999
$ws.Range("L123").Value = 89999
$ws.Range("N123").Value = -99799

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 4070.5715
$ws.Range("I11").Value = 5498.8
$ws.Range("K11").Value = 16496.4
$ws.Range("M11").Value = -16356.4
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
# Row 38
$ws.Range("H38").Value = 161
$ws.Range("J38").Value = 170.76923
$ws.Range("L38").Value = 512.30769
$ws.Range("N38").Value = -1206.30769
# Row 56
$ws.Range("H56").Value = 7228.675
$ws.Range("I56").Value = 7228.675
$ws.Range("K56").Value = 7228.675
$ws.Range("M56").Value = -6698.675
# Row 60
$ws.Range("H60").Value = 848.36365
$ws.Range("I60").Value = 237.125
$ws.Range("J60").Value = 2478.3333
$ws.Range("K60").Value = 711.375
$ws.Range("L60").Value = 7434.999899999999
$ws.Range("M60").Value = -460.375
$ws.Range("N60").Value = -7936.999899999999
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 262.44446
$ws.Range("I2").Value = 301.7857
$ws.Range("J2").Value = 124.75
$ws.Range("K2").Value = 301.7857
$ws.Range("L2").Value = 124.75
$ws.Range("M2").Value = -188.7857
$ws.Range("N2").Value = -350.75
# Row 55
$ws.Range("H55").Value = 30666.334
$ws.Range("J55").Value = 33000
$ws.Range("L55").Value = 33000
$ws.Range("N55").Value = -33654
# Row 99
$ws.Range("H99").Value = 8056.5
$ws.Range("I99").Value = 8056.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 8056.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -5810.5
$ws.Range("N99").ClearContents()
# Row 107
$ws.Range("H107").Value = 421.19232
$ws.Range("I107").Value = 405.14285
$ws.Range("J107").Value = 439.91666
$ws.Range("K107").Value = 405.14285
$ws.Range("L107").Value = 439.91666
$ws.Range("M107").Value = 1514.85715
$ws.Range("N107").Value = -4279.91666
# Row 113
$ws.Range("H113").Value = 7876.5
$ws.Range("I113").Value = 8738.143
$ws.Range("K113").Value = 8738.143
$ws.Range("M113").Value = -6568.143
# Row 122
$ws.Range("H122").Value = 4161.125
$ws.Range("I122").Value = 3186.9443
$ws.Range("K122").Value = 9560.832900000001
$ws.Range("M122").Value = -7110.832900000001
# Row 126
$ws.Range("H126").Value = 3683.6316
$ws.Range("I126").Value = 3076.077
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 9228.231
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -6758.231
$ws.Range("N126").Value = -19940

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1308.7778
$ws.Range("I82").Value = 980.3333
$ws.Range("J82").Value = 1965.6666
$ws.Range("K82").Value = 980.3333
$ws.Range("L82").Value = 1965.6666
$ws.Range("M82").Value = -619.3333
$ws.Range("N82").Value = -2687.6666
# Row 85
$ws.Range("H85").Value = 1308.7778
$ws.Range("I85").Value = 980.3333
$ws.Range("J85").Value = 1965.6666
$ws.Range("K85").Value = 980.3333
$ws.Range("L85").Value = 1965.6666
$ws.Range("M85").Value = 267.6667
$ws.Range("N85").Value = -4461.6666

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 10000
$ws.Range("J81").Value = 10000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -22122
# Row 84
$ws.Range("H84").Value = 10000
$ws.Range("J84").Value = 10000
$ws.Range("L84").Value = 100000
$ws.Range("N84").Value = -110608
# Row 100
$ws.Range("H100").Value = 969.8
$ws.Range("I100").Value = 969.8
$ws.Range("K100").Value = 1939.6
$ws.Range("M100").Value = -1398.6
# Row 107
$ws.Range("H107").Value = 2354.318
$ws.Range("I107").Value = 2138.4375
$ws.Range("J107").Value = 2930
$ws.Range("K107").Value = 6415.3125
$ws.Range("L107").Value = 8790
$ws.Range("M107").Value = -4495.3125
$ws.Range("N107").Value = -12630
# Row 132
$ws.Range("H132").Value = 3527.8667
$ws.Range("J132").Value = 4038.3333
$ws.Range("L132").Value = 12114.9999
$ws.Range("N132").Value = -17174.9999
# Row 136
$ws.Range("H136").Value = 1601.0769
$ws.Range("I136").Value = 1055.75
$ws.Range("K136").Value = 3167.25
$ws.Range("M136").Value = -617.25
